$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F ("District") to hold the new "Address" column,
# derived from the second line of column B ("NAMES and Address").
$ws.Columns.Item(6).Insert()
$ws.Cells.Item(2, 6).Value = "Address"

for ($row = 3; $row -le 42; $row++) {
    $nameAndAddress = $ws.Cells.Item($row, 2).Value2
    $lines = $nameAndAddress -split "`n"
    if ($lines.Length -ge 2) {
        $addrLine = $lines[1]
    } else {
        $addrLine = ""
    }

    $parts = $addrLine -split ","
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }
    if ($trimmed.Length -gt 0 -and $trimmed[$trimmed.Length - 1] -match "^Gadag\.?$") {
        $trimmed = $trimmed[0..($trimmed.Length - 2)]
    }
    $address = [string]::Join("", $trimmed)

    $ws.Cells.Item($row, 6).Value = $address
}
